$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K (shifts old K:M -> L:N, P:R -> Q:S, etc.)
$ws.Columns("K").Insert()

# New header cell in the inserted column
$ws.Range("K1").Value = "Customer defect name"
$ws.Range("K1").Font.Color = 255

# Match the width of the neighboring "Production Date" column (J) for the
# newly inserted column (K) so both render at the same width.
$ws.Columns("K").ColumnWidth = 28.17

# Data validation source ranges need to track the column shift (P:R -> Q:S)
$ws.Range("D2").Validation.Formula1 = "=`$Q`$1:`$Q`$2"
$ws.Range("E2").Validation.Formula1 = "=`$R`$1:`$R`$2"
$ws.Range("M2").Validation.Formula1 = "=`$S`$1:`$S`$2"

# Restore the selection shown in the sheet view
$null = $ws.Range("I7").Select()
